{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// The diff splits two single-run paragraphs into multiple <w:t> runs\n// separated by <w:br/> (a manual line break), by inserting a Word line\n// break (represented as \"\\u000b\" in Office.js Range/Body .text) at\n// specific points inside the existing run text \u2014 without otherwise\n// touching formatting.\n\nconst body = context.document.body;\n\n// 1) \"Crit\u00e9rio\" run: break the line right before \"M1 = 0,6*NR + 0,4* NP\".\n{\n  const oldText =\n    \"Alunos com m\u00e9dia final igual ou superior a 5,0 estar\u00e3o aprovados, desde que tenham frequ\u00eancia m\u00ednima de 70% (regimental). Alunos com m\u00e9dia inferior a 3,0 e/ou frequ\u00eancia inferior a 70% estar\u00e3o reprovados (regimental). Alunos com m\u00e9dia superior ou igual a 3,0 e inferior a 5,0 e que tenham frequ\u00eancia m\u00ednima de 70% ser\u00e3o submetidos ao per\u00edodo de recupera\u00e7\u00e3o (regimental). A m\u00e9dia (M1) ser\u00e1 calculada de acordo com a equa\u00e7\u00e3o abaixo:M1 = 0,6*NR + 0,4* NP\";\n  const newText =\n    \"Alunos com m\u00e9dia final igual ou superior a 5,0 estar\u00e3o aprovados, desde que tenham frequ\u00eancia m\u00ednima de 70% (regimental). Alunos com m\u00e9dia inferior a 3,0 e/ou frequ\u00eancia inferior a 70% estar\u00e3o reprovados (regimental). Alunos com m\u00e9dia superior ou igual a 3,0 e inferior a 5,0 e que tenham frequ\u00eancia m\u00ednima de 70% ser\u00e3o submetidos ao per\u00edodo de recupera\u00e7\u00e3o (regimental). A m\u00e9dia (M1) ser\u00e1 calculada de acordo com a equa\u00e7\u00e3o abaixo:\\u000bM1 = 0,6*NR + 0,4* NP\";\n\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  }\n}\n\n// 2) \"Bibliografia\" run: break the line before each numbered reference\n//    (2), 3), 4)) so each ends up on its own line.\n{\n  const oldText =\n    \"1)FOX, R.W.; PRITCHARD, P.J.; McDONALD, A.T. Introdu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos Fluidos. Ed. Gen LTC, 7 ed, Rio de Janeiro/RJ, 2010.2)\u00c7ENGEL, Y.A.; COMBALA, J.M. Mec\u00e2nica dos Fluidos: fundamentos e aplica\u00e7\u00f5es. McGraw-Hill Education (AMGH Editora Ltda),  Porto Alegre/ RS, 2007.3)COUPER, JR.; PENNEY, W.R.; FAIR, J.R.; WALAS, S.M. Chemical Process Equipment: Selection and Design. Amsterdam: Elsevier, 2005.4)TROPEA, C.; YARIN, A.L.; FOSS, J.F. Handbook of Experimental Fluid Mechanics. Ed Springer. Springer-Verlag Berlin Heidelberg. 2007\";\n  const newText =\n    \"1)FOX, R.W.; PRITCHARD, P.J.; McDONALD, A.T. Introdu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos Fluidos. Ed. Gen LTC, 7 ed, Rio de Janeiro/RJ, 2010.\\u000b2)\u00c7ENGEL, Y.A.; COMBALA, J.M. Mec\u00e2nica dos Fluidos: fundamentos e aplica\u00e7\u00f5es. McGraw-Hill Education (AMGH Editora Ltda),  Porto Alegre/ RS, 2007.\\u000b3)COUPER, JR.; PENNEY, W.R.; FAIR, J.R.; WALAS, S.M. Chemical Process Equipment: Selection and Design. Amsterdam: Elsevier, 2005.\\u000b4)TROPEA, C.; YARIN, A.L.; FOSS, J.F. Handbook of Experimental Fluid Mechanics. Ed Springer. Springer-Verlag Berlin Heidelberg. 2007\";\n\n  const results2 = body.search(oldText, { matchCase: true });\n  results2.load(\"items\");\n  await context.sync();\n\n  if (results2.items.length > 0) {\n    results2.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is pre-seeded as the target document.\n#\n# The diff splits two single-run paragraphs into multiple <w:t> runs\n# separated by <w:br/> (a manual line break) by inserting a line break\n# at specific points inside the existing run text, leaving formatting\n# untouched. We use Find/Replace with the \"^l\" (manual line break)\n# special code in the Replacement.Text, which is the standard Word way\n# to inject a <w:br/> without disturbing run formatting.\n\n$d = $word.ActiveDocument\n\nfunction Insert-LineBreakAfter($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $findText\n    $rng.Find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceOne = 1 (last arg selects replace mode);\n    # MatchCase = $true, Wrap = wdFindContinue (1), Replace = wdReplaceOne (1)\n    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n}\n\n# 1) \"Crit\u00e9rio\" run: break the line right before \"M1 = 0,6*NR + 0,4* NP\".\nInsert-LineBreakAfter \"equa\u00e7\u00e3o abaixo:M1 = 0,6*NR + 0,4* NP\" \"equa\u00e7\u00e3o abaixo:^lM1 = 0,6*NR + 0,4* NP\"\n\n# 2) \"Bibliografia\" run: break the line before each numbered reference\n#    (2), 3), 4)) so each ends up on its own line.\nInsert-LineBreakAfter \"Rio de Janeiro/RJ, 2010.2)\u00c7ENGEL\" \"Rio de Janeiro/RJ, 2010.^l2)\u00c7ENGEL\"\nInsert-LineBreakAfter \"Porto Alegre/ RS, 2007.3)COUPER\" \"Porto Alegre/ RS, 2007.^l3)COUPER\"\nInsert-LineBreakAfter \"Amsterdam: Elsevier, 2005.4)TROPEA\" \"Amsterdam: Elsevier, 2005.^l4)TROPEA\"\n"}
